$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a text value, forcing text storage even when the
# string looks like a number (e.g. "320.98"), then restore the cell style
# so no stray number-format style is left behind.
function Set-TextValue($cellRange, $value) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = "43.535.75"
$ws.Range("E2").Value = "  +2.43%  "

# Row 3
$ws.Range("D3").Value = "2.376.21"
$ws.Range("E3").Value = "  +6.67%  "

# Row 4
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
Set-TextValue $ws.Range("D5") "320.98"
$ws.Range("E5").Value = "  +8.17%  "

# Row 6
Set-TextValue $ws.Range("D6") "107.74"
$ws.Range("E6").Value = "  -3.45%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.639"
$ws.Range("E7").Value = "  +2.21%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("E9").Value = "  +4.32%  "

# Row 10
Set-TextValue $ws.Range("D10") "42.62"
$ws.Range("E10").Value = "  -3.93%  "

# Row 11
$ws.Range("E11").Value = "  +2.64%  "

# Row 12
Set-TextValue $ws.Range("D12") "8.73"
$ws.Range("E12").Value = "  -1.10%  "

# Row 13
$ws.Range("E13").Value = "  +3.33%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "16.69"
$ws.Range("E14").Value = "  +10.49%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D15") "0.106"
$ws.Range("E15").Value = "  +2.24%  "

# Row 16
$ws.Range("D16").Value = "2.734.76"
$ws.Range("E16").Value = "  +6.80%  "

# Row 17
$ws.Range("D17").Value = "2.413.94"
$ws.Range("E17").Value = "  +7.81%  "

# Row 18
$ws.Range("D18").Value = "43.547.69"
$ws.Range("E18").Value = "  +2.47%  "

# Row 19
$ws.Range("E19").Value = "  +3.39%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.27"
$ws.Range("E20").Value = "  -0.59%  "

# Row 21
Set-TextValue $ws.Range("D21") "75.65"
$ws.Range("E21").Value = "  +3.69%  "

# Row 22
Set-TextValue $ws.Range("D22") "3.52"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
Set-TextValue $ws.Range("D23") "263.84"
$ws.Range("E23").Value = "  +14.81%  "

# Row 24
$ws.Range("E24").Value = "  +6.08%  "

# Row 25
Set-TextValue $ws.Range("D25") "9.19"
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
Set-TextValue $ws.Range("D26") "12.03"
$ws.Range("E26").Value = "  +3.62%  "

# Row 27
$ws.Range("E27").Value = "  +0.02%  "

# Row 28
Set-TextValue $ws.Range("D28") "39.05"
$ws.Range("E28").Value = "  +1.64%  "

# Row 29
Set-TextValue $ws.Range("D29") "22.92"
$ws.Range("E29").Value = "  +8.94%  "

# Row 30
$ws.Range("E30").Value = "  +0.42%  "

# Row 31
$ws.Range("E31").Value = "  -0.81%  "

# Row 32
Set-TextValue $ws.Range("D32") "174.43"
$ws.Range("E32").Value = "  +0.67%  "

# Row 33
$ws.Range("E33").Value = "  +3.18%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.96"
$ws.Range("E34").Value = "  +3.24%  "

# Row 35
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D35") "4.97"
$ws.Range("E35").Value = "  -3.44%  "

# Row 36
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D36") "0.131"
$ws.Range("E36").Value = "  +3.92%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D37") "4.11"
$ws.Range("E37").Value = "  -3.81%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.0374"
$ws.Range("E38").Value = "  -0.80%  "

# Row 39
$ws.Range("E39").Value = "  +1.62%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.84"
$ws.Range("E40").Value = "  +17.95%  "

# Row 41
$ws.Range("E41").Value = "  +13.65%  "

# Row 42
Set-TextValue $ws.Range("D42") "71.88"
$ws.Range("E42").Value = "  -0.32%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.231"
$ws.Range("E43").Value = "  -1.70%  "

# Row 44
$ws.Range("E44").Value = "  +0.11%  "

# Row 45
Set-TextValue $ws.Range("D45") "12.64"
$ws.Range("E45").Value = "  -0.83%  "

# Row 46
Set-TextValue $ws.Range("D46") "5.63"
$ws.Range("E46").Value = "  +3.62%  "

# Row 47
Set-TextValue $ws.Range("D47") "113.02"
$ws.Range("E47").Value = "  +9.69%  "

# Row 48
$ws.Range("E48").Value = "  +8.54%  "

# Row 49
$ws.Range("E49").Value = "  -0.50%  "

# Row 50
$ws.Range("E50").Value = "  +3.16%  "

# Row 51
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D51") "71.54"
$ws.Range("E51").Value = "  +5.05%  "
